$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("A3").Value = "on Little Master"
$ws.Range("A3").Select()
